# Gravity Boii Runner-Estimation.xlsx update
# "Particle Effects all set up and ready. Fixing Dash and Shift Effects"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Particle Effects sub-items (rows 39-41) ---
# Row 39: "Dash flames" -> "Dash Sprites behind Trail", now has a Time Taken value (Done)
$ws.Range("B39").Value = "Dash Sprites behind Trail"
$ws.Range("D39").Value = 2

# Row 40: "Shift wind." -> "Shift gust of wind trail", now has a Time Taken value (Done)
$ws.Range("B40").Value = "Shift gust of wind trail"
$ws.Range("D40").Value = 2

# Row 41: "Sparks for gravity pool" text unchanged, but now has a Time Taken value (Done)
$ws.Range("D41").Value = 2

# Recolor rows 39-41 (B:D) from the "not done" yellow to the "done" green,
# matching the rest of the completed sub-items (e.g. rows 36-37).
$ws.Range("B39:D40").Interior.Color = 5296274
$ws.Range("B41:D41").Interior.Color = 5296274

# Row 38 status note updated to reflect current progress on Particle Effects
$ws.Range("E38").Value = "Dash and Shift Particles set up but not showing"

# --- Update selection / scroll position ---
$ws.Range("E42:E46").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
